# Generate Report for Archive
#
# 1. Status text moved from "Ready for handoff" to "In Translation" on every
#    sheet that surfaces it (Overview!E2/F2 plus the per-language status
#    column C2 on each language sheet).
# 2. The status columns got narrower to match the new (shorter) text:
#    Overview columns E & F, and column C on each language sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears ------------------
foreach ($sheet in $wb.Worksheets) {
    [void]$sheet.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Narrow the status columns to the new auto-fit width -----------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
